$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column D, shifting old D:K to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# New columns should carry the same number formatting as the (now-shifted)
# original data columns; copy formats only (not values) from F to D:E.
# (Rows 36:37 and 78:79 are blank separator/header rows with no data
# columns at all, so they are intentionally excluded.)
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the column width of the neighboring (formerly D) data column.
$colWidth = $ws.Columns("F").ColumnWidth
$ws.Columns("D").ColumnWidth = $colWidth
$ws.Columns("E").ColumnWidth = $colWidth

# Populate the two new quarter columns (D = latest quarter, E = prior quarter)
# with the reported figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1038000
$ws.Range("E8").Value = 1080000
$ws.Range("D9").Value = 897000
$ws.Range("E9").Value = 921000
$ws.Range("D10").Value = 141000
$ws.Range("E10").Value = 159000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = -22000
$ws.Range("E14").Value = 20000
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 940000
$ws.Range("E17").Value = 1023000
$ws.Range("D18").Value = 98000
$ws.Range("E18").Value = 57000
$ws.Range("D20").Value = 15000
$ws.Range("E20").Value = -6000
$ws.Range("D21").Value = 135000
$ws.Range("E21").Value = 71000
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 113000
$ws.Range("E23").Value = 51000
$ws.Range("D24").Value = 28000
$ws.Range("E24").Value = 6000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 85000
$ws.Range("E26").Value = 45000
$ws.Range("D27").Value = 83000
$ws.Range("E27").Value = 44000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 7000
$ws.Range("E29").Value = -12000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -15000
$ws.Range("E32").Value = 6000
$ws.Range("D33").Value = 90000
$ws.Range("E33").Value = 32000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 90000
$ws.Range("E35").Value = 32000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 127000
$ws.Range("E41").Value = 115000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 547000
$ws.Range("E43").Value = 604000
$ws.Range("D44").Value = 522000
$ws.Range("E44").Value = 477000
$ws.Range("D45").Value = 44000
$ws.Range("E45").Value = 30000
$ws.Range("D46").Value = 1240000
$ws.Range("E46").Value = 1226000
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = 178000
$ws.Range("D48").Value = 481000
$ws.Range("E48").Value = 483000
$ws.Range("D49").Value = 418000
$ws.Range("E49").Value = 465000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 594000
$ws.Range("E52").Value = 374000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2733000
$ws.Range("E54").Value = 2726000
$ws.Range("D57").Value = 680000
$ws.Range("E57").Value = 700000
$ws.Range("D58").Value = 172000
$ws.Range("E58").Value = 94000
$ws.Range("D59").Value = 300000
$ws.Range("E59").Value = 290000
$ws.Range("D60").Value = 1152000
$ws.Range("E60").Value = 1084000
$ws.Range("D61").Value = 731000
$ws.Range("E61").Value = 730000
$ws.Range("D62").Value = 488000
$ws.Range("E62").Value = 595000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 2404000
$ws.Range("E66").Value = 2439000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 290000
$ws.Range("E72").Value = 200000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 329000
$ws.Range("E76").Value = 287000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 90000
$ws.Range("E81").Value = 32000
$ws.Range("D83").Value = 22000
$ws.Range("E83").Value = 20000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 11000
$ws.Range("E89").Value = 60000
$ws.Range("D91").Value = -23000
$ws.Range("E91").Value = -52000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -27000
$ws.Range("E94").Value = -51000
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 27000
$ws.Range("E100").Value = 7000
$ws.Range("D101").Value = 1000
$ws.Range("E101").Value = -1000
$ws.Range("D102").Value = 12000
$ws.Range("E102").Value = 15000

